$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers/values for D2:G5 (column D = "U" strings, E/F/G = numbers)
$ws.Range("D2").Value = "U"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("D3").Value = "U"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2

$ws.Range("D4").Value = "U"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2

$ws.Range("D5").Value = "U"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3

# Update selection to match the new data range
$ws.Range("D2:G5").Select()
